$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-01-17 Wednesday" "2024-01-18 Thursday"

Replace-Text "81×74=5994" "46×57=2622"
Replace-Text "27×11=297" "50×17=850"
Replace-Text "32×38=1216" "74×76=5624"
Replace-Text "51×94=4794" "27×87=2349"
Replace-Text "75×66=4950" "38×22=836"

Replace-Text "81×91=7371" "30×16=480"
Replace-Text "29×54=1566" "34×44=1496"
Replace-Text "71×62=4402" "87×31=2697"
Replace-Text "39×26=1014" "67×51=3417"
Replace-Text "31×63=1953" "23×72=1656"

Replace-Text "73×28=2044" "20×14=280"
Replace-Text "48×42=2016" "18×40=720"
Replace-Text "63×39=2457" "96×93=8928"
Replace-Text "67×65=4355" "44×54=2376"
Replace-Text "14×62=868" "75×60=4500"

Replace-Text "88×80=7040" "15×66=990"
Replace-Text "67×92=6164" "66×45=2970"
Replace-Text "54×51=2754" "38×89=3382"
Replace-Text "35×63=2205" "58×25=1450"
Replace-Text "64×24=1536" "39×18=702"

Replace-Text "23×81=1863" "39×64=2496"
Replace-Text "34×31=1054" "14×11=154"
Replace-Text "27×57=1539" "98×34=3332"
Replace-Text "18×39=702" "63×50=3150"
Replace-Text "57×17=969" "85×47=3995"
